$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.NumberFormat = "General"
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.110.04"
Set-TextValue $ws.Range("E2") "  +1.28%  "
Set-TextValue $ws.Range("D3") "3.553.94"
Set-TextValue $ws.Range("E3") "  +1.91%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "619.21"
Set-TextValue $ws.Range("E5") "  +2.16%  "
Set-TextValue $ws.Range("D6") "154.72"
Set-TextValue $ws.Range("E6") "  +4.18%  "
Set-TextValue $ws.Range("D7") "3.553.43"
Set-TextValue $ws.Range("E7") "  +1.91%  "
Set-TextValue $ws.Range("E8") "  +0.02%  "
Set-TextValue $ws.Range("E9") "  +1.87%  "
Set-TextValue $ws.Range("E10") "  +5.37%  "
Set-TextValue $ws.Range("D11") "7.32"
Set-TextValue $ws.Range("E11") "  +5.33%  "
Set-TextValue $ws.Range("E12") "  +3.35%  "
Set-TextValue $ws.Range("B13") "Avalanche"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D13") "33.17"
Set-TextValue $ws.Range("E13") "  +4.92%  "
Set-TextValue $ws.Range("B14") "ShibaInu"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D14") "0.0000221"
Set-TextValue $ws.Range("E14") "  +1.67%  "
Set-TextValue $ws.Range("D15") "4.154.85"
Set-TextValue $ws.Range("E15") "  +1.95%  "
Set-TextValue $ws.Range("D16") "3.560.81"
Set-TextValue $ws.Range("E16") "  +2.03%  "
Set-TextValue $ws.Range("D17") "68.058.73"
Set-TextValue $ws.Range("E17") "  +1.35%  "
Set-TextValue $ws.Range("E18") "  +0.01%  "
Set-TextValue $ws.Range("E19") "  +6.06%  "
Set-TextValue $ws.Range("E20") "  +5.69%  "
Set-TextValue $ws.Range("D21") "9.95"
Set-TextValue $ws.Range("E21") "  +10.26%  "
Set-TextValue $ws.Range("D22") "454.40"
Set-TextValue $ws.Range("E22") "  +1.49%  "
Set-TextValue $ws.Range("D23") "0.641"
Set-TextValue $ws.Range("E23") "  +2.86%  "
Set-TextValue $ws.Range("D24") "78.26"
Set-TextValue $ws.Range("E24") "  +1.40%  "
Set-TextValue $ws.Range("E25") "  +4.02%  "
Set-TextValue $ws.Range("B26") "WrappedeETH"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D26") "3.694.64"
Set-TextValue $ws.Range("E26") "  +1.89%  "
Set-TextValue $ws.Range("B27") "PEPE"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D27") "0.0000129"
Set-TextValue $ws.Range("E27") "  +1.31%  "
Set-TextValue $ws.Range("D28") "0.999"
Set-TextValue $ws.Range("E28") "  -0.11%  "
Set-TextValue $ws.Range("D29") "9.07"
Set-TextValue $ws.Range("E29") "  +9.19%  "
Set-TextValue $ws.Range("E30") "  +2.81%  "
Set-TextValue $ws.Range("E31") "  +7.42%  "
Set-TextValue $ws.Range("D32") "0.171"
Set-TextValue $ws.Range("E32") "  +4.28%  "
Set-TextValue $ws.Range("D33") "1.00"
Set-TextValue $ws.Range("E33") "  +0.04%  "
Set-TextValue $ws.Range("D34") "6.41"
Set-TextValue $ws.Range("E34") "  +4.56%  "
Set-TextValue $ws.Range("D35") "26.04"
Set-TextValue $ws.Range("E35") "  +1.29%  "
Set-TextValue $ws.Range("E36") "  +3.63%  "
Set-TextValue $ws.Range("D37") "3.545.98"
Set-TextValue $ws.Range("E37") "  +1.94%  "
Set-TextValue $ws.Range("D38") "8.26"
Set-TextValue $ws.Range("E38") "  +3.36%  "
Set-TextValue $ws.Range("D39") "2.36"
Set-TextValue $ws.Range("E39") "  +7.29%  "
Set-TextValue $ws.Range("E40") "  -0.02%  "
Set-TextValue $ws.Range("D41") "178.52"
Set-TextValue $ws.Range("E41") "  +3.76%  "
Set-TextValue $ws.Range("E42") "  +5.32%  "
Set-TextValue $ws.Range("D43") "1.00"
Set-TextValue $ws.Range("E43") "  +0.08%  "
Set-TextValue $ws.Range("D44") "5.60"
Set-TextValue $ws.Range("E44") "  +3.33%  "
Set-TextValue $ws.Range("D45") "31.00"
Set-TextValue $ws.Range("E45") "  +16.02%  "
Set-TextValue $ws.Range("E46") "  +1.39%  "
Set-TextValue $ws.Range("D47") "46.47"
Set-TextValue $ws.Range("E47") "  +2.42%  "
Set-TextValue $ws.Range("D48") "1.34"
Set-TextValue $ws.Range("E48") "  +7.12%  "
Set-TextValue $ws.Range("E49") "  +3.45%  "
Set-TextValue $ws.Range("E50") "  +3.13%  "
Set-TextValue $ws.Range("E51") "  +6.71%  "
